# Auto-generated edit script: refresh cryptos list values (cell-by-cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/D hold text that can look numeric/date-like (e.g. "26.250.54",
# "0.5244", "8.029"). Assigning those through .Value would let Excel auto-
# coerce/parse them, so we force plain text with a leading apostrophe and then
# restore the default "Normal" style so no stray quote-prefix formatting is left
# behind (matching the original unstyled inline-string cells).
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.250.54'
$ws.Range('E2').Value = '  -0.41%  '
Set-TextValue $ws.Range('D3') '1.658.64'
$ws.Range('E4').Value = '  -0.73%  '
Set-TextValue $ws.Range('D5') '219.22'
$ws.Range('E5').Value = '  -0.10%  '
Set-TextValue $ws.Range('D6') '0.5244'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  -0.59%  '
Set-TextValue $ws.Range('D8') '0.2646'
$ws.Range('E8').Value = '  -0.49%  '
Set-TextValue $ws.Range('D9') '0.06322'
$ws.Range('E9').Value = '  -1.13%  '
Set-TextValue $ws.Range('D10') '20.69'
$ws.Range('E10').Value = '  -0.90%  '
Set-TextValue $ws.Range('D11') '0.07787'
$ws.Range('E11').Value = '  -0.88%  '
Set-TextValue $ws.Range('D12') '4.514'
$ws.Range('E12').Value = '  -1.08%  '
Set-TextValue $ws.Range('D13') '1.563.74'
$ws.Range('E13').Value = '  -6.33%  '
Set-TextValue $ws.Range('D14') '1.887.74'
$ws.Range('E14').Value = '  -0.44%  '
Set-TextValue $ws.Range('D15') '0.5634'
$ws.Range('E15').Value = '  +1.75%  '
Set-TextValue $ws.Range('D16') '0.0₅8069'
$ws.Range('E16').Value = '  -1.32%  '
Set-TextValue $ws.Range('D17') '65.31'
$ws.Range('E17').Value = '  -1.00%  '
Set-TextValue $ws.Range('D18') '26.246.86'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('E19').Value = '  -0.71%  '
Set-TextValue $ws.Range('D20') '4.727'
$ws.Range('E20').Value = '  +1.16%  '
Set-TextValue $ws.Range('D21') '194.34'
$ws.Range('E21').Value = '  -0.32%  '
Set-TextValue $ws.Range('D22') '10.25'
$ws.Range('E22').Value = '  -0.08%  '
Set-TextValue $ws.Range('D23') '6.036'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.73%  '
Set-TextValue $ws.Range('D25') '145.37'
$ws.Range('E25').Value = '  -0.88%  '
Set-TextValue $ws.Range('D26') '0.1211'
$ws.Range('E26').Value = '  -1.14%  '
Set-TextValue $ws.Range('D27') '7.231'
$ws.Range('E27').Value = '  -0.11%  '
Set-TextValue $ws.Range('D28') '16.04'
$ws.Range('E28').Value = '  -1.04%  '
Set-TextValue $ws.Range('D29') '1.504'
$ws.Range('E29').Value = '  +0.12%  '
Set-TextValue $ws.Range('D30') '0.05642'
$ws.Range('E30').Value = '  -3.68%  '
Set-TextValue $ws.Range('D31') '1.278'
$ws.Range('E31').Value = '  -0.62%  '
Set-TextValue $ws.Range('D32') '3.489'
$ws.Range('E32').Value = '  -2.52%  '
Set-TextValue $ws.Range('D33') '3.376'
$ws.Range('E33').Value = '  +2.70%  '
Set-TextValue $ws.Range('D35') '2.807'
$ws.Range('E35').Value = '  -0.88%  '
Set-TextValue $ws.Range('D36') '0.9468'
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  -1.20%  '
Set-TextValue $ws.Range('D39') '0.01606'
$ws.Range('E39').Value = '  +0.04%  '
Set-TextValue $ws.Range('D40') '5.982'
$ws.Range('E40').Value = '  +2.06%  '
Set-TextValue $ws.Range('D41') '2.570'
$ws.Range('E41').Value = '  -0.74%  '
Set-TextValue $ws.Range('D42') '1.047.60'
$ws.Range('E42').Value = '  -2.71%  '
Set-TextValue $ws.Range('D43') '0.8467'
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('E44').Value = '  -0.75%  '
Set-TextValue $ws.Range('D45') '102.87'
$ws.Range('E45').Value = '  -1.38%  '
Set-TextValue $ws.Range('D46') '1.799.08'
$ws.Range('E46').Value = '  -0.38%  '
Set-TextValue $ws.Range('D47') '58.47'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  +0.09%  '
Set-TextValue $ws.Range('D49') '1.002'
$ws.Range('E49').Value = '  -1.19%  '
Set-TextValue $ws.Range('D50') '0.05331'
$ws.Range('E50').Value = '  +3.13%  '
Set-TextValue $ws.Range('B51') 'EnergySwap'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '8.029'
$ws.Range('E51').Value = '  +0.11%  '
